# Insert a new data row at row 59 (pushing the existing rows 59-142 down to
# 60-143) and populate it with the new Acelga price-report entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Insert()

$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(59, 3).Value = 'Los Lagos'
$ws.Cells.Item(59, 4).Value = 44580
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112009
$ws.Cells.Item(59, 7).Value = 'Acelga'
$ws.Cells.Item(59, 8).Value = 'Sin especificar'
$ws.Cells.Item(59, 9).Value = 'Primera'
$ws.Cells.Item(59, 10).Value = 20
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 10000
$ws.Cells.Item(59, 13).Value = 10000
$ws.Cells.Item(59, 14).Value = '$/docena de atados (12 kilos)'
$ws.Cells.Item(59, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(59, 16).Value = 833
$ws.Cells.Item(59, 17).Value = 12
$ws.Cells.Item(59, 18).Value = 'Hortaliza'
